$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.998.20"
$ws.Range("E2").Value = "  -1.13%  "
$ws.Range("D3").Value = "2.662.75"
$ws.Range("E3").Value = "  +1.21%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'598.15"
$ws.Range("E5").Value = "  -1.15%  "
$ws.Range("D6").Value = "'174.55"
$ws.Range("E6").Value = "  -2.36%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "'0.523"
$ws.Range("E8").Value = "  -0.69%  "
$ws.Range("D9").Value = "2.663.92"
$ws.Range("E9").Value = "  +1.31%  "
$ws.Range("E10").Value = "  -1.46%  "
$ws.Range("E11").Value = "  +2.44%  "
$ws.Range("E12").Value = "  +0.75%  "
$ws.Range("E13").Value = "  -1.27%  "
$ws.Range("D14").Value = "3.149.24"
$ws.Range("E15").Value = "  -2.66%  "
$ws.Range("D16").Value = "71.874.85"
$ws.Range("E16").Value = "  -1.03%  "
$ws.Range("D17").Value = "'26.23"
$ws.Range("E17").Value = "  -1.70%  "
$ws.Range("D18").Value = "2.664.58"
$ws.Range("E18").Value = "  +1.28%  "
$ws.Range("D19").Value = "'12.26"
$ws.Range("E19").Value = "  +6.21%  "
$ws.Range("D20").Value = "'8.19"
$ws.Range("E20").Value = "  +1.96%  "
$ws.Range("D21").Value = "'370.95"
$ws.Range("E21").Value = "  -3.51%  "
$ws.Range("D22").Value = "'4.16"
$ws.Range("E22").Value = "  -0.46%  "
$ws.Range("D23").Value = "'2.04"
$ws.Range("E23").Value = "  +1.49%  "
$ws.Range("D24").Value = "'72.10"
$ws.Range("E24").Value = "  -1.41%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("E26").Value = "  -1.44%  "
$ws.Range("E27").Value = "  -1.39%  "
$ws.Range("D28").Value = "2.801.47"
$ws.Range("E28").Value = "  +1.43%  "
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("D30").Value = "0.0₃0967"
$ws.Range("E30").Value = "  +0.84%  "
$ws.Range("D31").Value = "'8.04"
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("D32").Value = "'500.28"
$ws.Range("E32").Value = "  -5.68%  "
$ws.Range("E33").Value = "  -2.28%  "
$ws.Range("E34").Value = "  -0.49%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.21%  "
$ws.Range("D36").Value = "'163.05"
$ws.Range("E36").Value = "  -0.74%  "
$ws.Range("D37").Value = "'19.53"
$ws.Range("E37").Value = "  +0.94%  "
$ws.Range("D38").Value = "'19.05"
$ws.Range("E38").Value = "  -0.31%  "
$ws.Range("D39").Value = "'0.111"
$ws.Range("E39").Value = "  +0.36%  "
$ws.Range("E40").Value = "  -2.05%  "
$ws.Range("E41").Value = "  -3.60%  "
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("D43").Value = "'4.98"
$ws.Range("E43").Value = "  -1.64%  "
$ws.Range("E44").Value = "  -2.02%  "
$ws.Range("E45").Value = "  -0.28%  "
$ws.Range("D46").Value = "'156.21"
$ws.Range("E46").Value = "  +3.31%  "
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("D48").Value = "'0.559"
$ws.Range("E48").Value = "  +3.18%  "
$ws.Range("D49").Value = "'3.72"
$ws.Range("E49").Value = "  +0.82%  "
$ws.Range("E50").Value = "  +2.00%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "'0.606"
$ws.Range("E51").Value = "  +0.91%  "
